# Auto-generated edit script: applies the 2022-09-07 data refresh
# to the CTA index crime YTD workbook (167 cell updates + 1 new cell)
# across 27 worksheets (Citywide Totals, By Neighborhood, and 25
# individual neighborhood sheets).

$wb = $excel.ActiveWorkbook

# --- Citywide Totals (sheet1) ---
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("C2").Value = 49   # was 48
$ws.Range("I2").Value = 83   # was 82
$ws.Range("B3").Value = 61   # was 60
$ws.Range("F3").Value = 94   # was 93
$ws.Range("G3").Value = 91   # was 89
$ws.Range("H3").Value = 94   # was 93
$ws.Range("I3").Value = 142   # was 140
$ws.Range("D9").Value = 306   # was 305
$ws.Range("E9").Value = 310   # was 308
$ws.Range("F9").Value = 391   # was 390
$ws.Range("G9").Value = 352   # was 350
$ws.Range("H9").Value = 319   # was 316
$ws.Range("I9").Value = 377   # was 376
$ws.Range("B10").Value = 904   # was 898
$ws.Range("C10").Value = 1098   # was 1091
$ws.Range("D10").Value = 1261   # was 1256
$ws.Range("E10").Value = 1556   # was 1550
$ws.Range("F10").Value = 1593   # was 1587
$ws.Range("G10").Value = 763   # was 762
$ws.Range("H10").Value = 391   # was 387
$ws.Range("I10").Value = 621   # was 619
$ws.Range("B11").Value = 1280   # was 1273
$ws.Range("C11").Value = 1562   # was 1554
$ws.Range("D11").Value = 1737   # was 1731
$ws.Range("E11").Value = 2025   # was 2017
$ws.Range("F11").Value = 2148   # was 2140
$ws.Range("G11").Value = 1277   # was 1272
$ws.Range("H11").Value = 897   # was 889
$ws.Range("I11").Value = 1248   # was 1242

# --- Chinatown (sheet10) ---
$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("I3").Value = 5   # was 4
$ws.Range("I9").Value = 24   # was 23

# --- Garfield Park (sheet12) ---
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("F8").Value = 101   # was 99
$ws.Range("H8").Value = 22   # was 21
$ws.Range("F9").Value = 153   # was 151
$ws.Range("H9").Value = 61   # was 60

# --- Chatham (sheet13) ---
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("E7").Value = 25   # was 24
$ws.Range("E8").Value = 38   # was 37

# --- Grand Crossing (sheet14) ---
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("B3").Value = 3   # was 2
$ws.Range("F3").Value = 15   # was 14
$ws.Range("H7").Value = 18   # was 16
$ws.Range("I7").Value = 14   # was 13
$ws.Range("H8").Value = 16   # was 15
$ws.Range("B9").Value = 45   # was 44
$ws.Range("F9").Value = 67   # was 66
$ws.Range("H9").Value = 46   # was 43
$ws.Range("I9").Value = 59   # was 58

# --- Loop (sheet15) ---
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("B8").Value = 123   # was 121
$ws.Range("C8").Value = 211   # was 209
$ws.Range("D8").Value = 368   # was 366
$ws.Range("E8").Value = 451   # was 448
$ws.Range("F8").Value = 425   # was 422
$ws.Range("I8").Value = 147   # was 146
$ws.Range("B9").Value = 159   # was 157
$ws.Range("C9").Value = 252   # was 250
$ws.Range("D9").Value = 426   # was 424
$ws.Range("E9").Value = 510   # was 507
$ws.Range("F9").Value = 480   # was 477
$ws.Range("I9").Value = 246   # was 245

# --- Old Town (sheet17) ---
$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("E5").Value = 9   # was 8
$ws.Range("E7").Value = 49   # was 48

# --- Little Italy, UIC (sheet18) ---
$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("D5").Value = 8   # was 7
$ws.Range("D7").Value = 24   # was 23

# --- By Neighborhood (sheet2) ---
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("B7").Value = 37   # was 36
$ws.Range("G8").Value = 75   # was 74
$ws.Range("E19").Value = 38   # was 37
$ws.Range("I21").Value = 24   # was 23
$ws.Range("E22").Value = 14   # was 13
$ws.Range("C23").Value = 15   # was 13
$ws.Range("B28").Value = 73   # was 72
$ws.Range("E28").Value = 69   # was 68
$ws.Range("G28").Value = 62   # was 61
$ws.Range("F32").Value = 153   # was 151
$ws.Range("H32").Value = 61   # was 60
$ws.Range("B36").Value = 45   # was 44
$ws.Range("F36").Value = 67   # was 66
$ws.Range("H36").Value = 46   # was 43
$ws.Range("I36").Value = 59   # was 58
$ws.Range("H43").Value = 5   # was 4
$ws.Range("G45").Value = 24   # was 23
$ws.Range("I47").Value = 41   # was 40
$ws.Range("D50").Value = 24   # was 23
$ws.Range("B53").Value = 159   # was 157
$ws.Range("C53").Value = 252   # was 250
$ws.Range("D53").Value = 426   # was 424
$ws.Range("E53").Value = 510   # was 507
$ws.Range("F53").Value = 480   # was 477
$ws.Range("I53").Value = 246   # was 245
$ws.Range("B63").Value = 8   # was 7
$ws.Range("E70").Value = 49   # was 48
$ws.Range("D72").Value = 8   # was 7
$ws.Range("E74").Value = 63   # was 62
$ws.Range("G74").Value = 29   # was 28
$ws.Range("D76").Value = 43   # was 42
$ws.Range("G76").Value = 33   # was 32
$ws.Range("F77").Value = 44   # was 43
$ws.Range("H77").Value = 41   # was 39
$ws.Range("I77").Value = 61   # was 59
$ws.Range("B78").Value = 30   # was 29
$ws.Range("C87").Value = 32   # was 31
$ws.Range("F87").Value = 22   # was 21
$ws.Range("H87").Value = 21   # was 20
$ws.Range("C92").Value = 15   # was 14
$ws.Range("C93").Value = 7   # was 6
$ws.Range("C97").Value = 16   # was 15
$ws.Range("D97").Value = 17   # was 16
$ws.Range("B99").Value = 1280   # was 1273
$ws.Range("C99").Value = 1562   # was 1554
$ws.Range("D99").Value = 1737   # was 1731
$ws.Range("E99").Value = 2025   # was 2017
$ws.Range("F99").Value = 2148   # was 2140
$ws.Range("G99").Value = 1277   # was 1272
$ws.Range("H99").Value = 897   # was 889
$ws.Range("I99").Value = 1248   # was 1242

# --- Uptown (sheet24) ---
$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("C2").Value = 3   # was 2
$ws.Range("F7").Value = 5   # was 4
$ws.Range("H8").Value = 9   # was 8
$ws.Range("C9").Value = 32   # was 31
$ws.Range("F9").Value = 22   # was 21
$ws.Range("H9").Value = 21   # was 20

# --- Rush & Division (sheet25) ---
$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("B5").Value = 26   # was 25
$ws.Range("B6").Value = 30   # was 29

# --- Englewood (sheet26) ---
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("G3").Value = 4   # was 3
$ws.Range("B8").Value = 45   # was 44
$ws.Range("E8").Value = 46   # was 45
$ws.Range("B9").Value = 73   # was 72
$ws.Range("E9").Value = 69   # was 68
$ws.Range("G9").Value = 62   # was 61

# --- Lake View (sheet27) ---
$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("I2").Value = 2   # was 1
$ws.Range("I8").Value = 41   # was 40

# --- Jefferson Park (sheet28) ---
$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("G6").Value = 20   # was 19
$ws.Range("G7").Value = 24   # was 23

# --- Rogers Park (sheet3) ---
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("G6").Value = 6   # was 5
$ws.Range("D7").Value = 33   # was 32
$ws.Range("D8").Value = 43   # was 42
$ws.Range("G8").Value = 33   # was 32

# --- River North (sheet31) ---
$ws = $wb.Worksheets.Item("River North")
$ws.Range("G5").Value = 11   # was 10
$ws.Range("E6").Value = 59   # was 58
$ws.Range("E7").Value = 63   # was 62
$ws.Range("G7").Value = 29   # was 28

# --- West Loop (sheet33) ---
$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("C8").Value = 13   # was 12
$ws.Range("C9").Value = 15   # was 14

# --- Woodlawn (sheet37) ---
$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("C6").Value = 6   # was 5
$ws.Range("D6").Value = 10   # was 9
$ws.Range("C7").Value = 16   # was 15
$ws.Range("D7").Value = 17   # was 16

# --- Douglas (sheet46) ---
$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("C5").Value = 10   # was 8
$ws.Range("C6").Value = 15   # was 13

# --- West Pullman (sheet47) ---
$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("C5").Value = 4   # was 3
$ws.Range("C6").Value = 7   # was 6

# --- Roseland (sheet5) ---
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("H3").Value = 5   # was 4
$ws.Range("F8").Value = 31   # was 30
$ws.Range("H8").Value = 20   # was 19
$ws.Range("I8").Value = 28   # was 26
$ws.Range("F9").Value = 44   # was 43
$ws.Range("H9").Value = 41   # was 39
$ws.Range("I9").Value = 61   # was 59

# --- Clearing (sheet54) ---
$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("E7").Value = 14   # was 13
$ws.Range("E5").Value = 1   # new cell

# --- New City (sheet58) ---
$ws = $wb.Worksheets.Item("New City")
$ws.Range("B5").Value = 4   # was 3
$ws.Range("B6").Value = 8   # was 7

# --- Auburn Gresham (sheet6) ---
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("B6").Value = 23   # was 22
$ws.Range("B7").Value = 37   # was 36

# --- Printers Row (sheet62) ---
$ws = $wb.Worksheets.Item("Printers Row")
$ws.Range("D5").Value = 4   # was 3
$ws.Range("D6").Value = 8   # was 7

# --- Austin (sheet8) ---
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("G3").Value = 5   # was 4
$ws.Range("G8").Value = 75   # was 74

# --- Irving Park (sheet82) ---
$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("H5").Value = 3   # was 2
$ws.Range("H7").Value = 5   # was 4
